# "Now generating a lattice graph for factors of production for each country."
#
# The CNData tab holds the indexed series re-exported for plotting in R.
# To support a per-country lattice plot, a "Year" column is inserted at
# the front (pulled straight from the China Workbook tab) and two trailing
# columns are added that tag every row with the indexing base ("iU" / "NA")
# and the country code ("Country" / "CN") so multiple countries' data can
# later be stacked into one long data frame.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CNData")

# Shift the existing iYear/iGDP/iLabor/iCapStk/iQ/iX columns one to the
# right, making room for the new leading "Year" column.
$ws.Columns.Item(1).Insert()

# Header row.
$ws.Range("A1").Formula = "='China Workbook'!A11"

# Year values, row by row (China Workbook rows 12-32 -> CNData rows 2-22).
for ($offset = 0; $offset -le 20; $offset++) {
    $sourceRow = 12 + $offset
    $targetRow = 2 + $offset
    $ws.Range("A" + $targetRow).Formula = "='China Workbook'!A" + $sourceRow
}

# New trailing columns: iU / Country headers, NA / CN values for every
# data row.
$ws.Range("H1").Value = "iU"
$ws.Range("I1").Value = "Country"

for ($row = 2; $row -le 22; $row++) {
    $ws.Range("H" + $row).Value = "NA"
    $ws.Range("I" + $row).Value = "CN"
}

# Keep the tab's selection sane after widening the used range.
$ws.Range("A1:I1048576").Select()
